$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Stato Attuale" (current fleet assignment state) ---
$ws1 = $wb.Worksheets.Item("Stato Attuale")

# Capture the operator currently assigned to targa GH228TC (row 3) before changing it,
# so it can be logged into the history sheet as the "previous operator".
$targa = $ws1.Range("A3").Value()
$oldOperatore = $ws1.Range("B3").Value()
$newOperatore = "MEZZANOTTE SABRINA"
$changeDate = "2026-01-05"

# Update the current state: new operator assigned, and record the assignment date.
$ws1.Range("B3").Value = $newOperatore

$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = $changeDate
$ws1.Range("C3").Style = "Normal"

# --- Sheet 2: "Storico Passaggi" (history log of operator changes) ---
$ws2 = $wb.Worksheets.Item("Storico Passaggi")

$ws2.Range("A3").Value = $targa
$ws2.Range("B3").Value = $oldOperatore
$ws2.Range("C3").Value = $newOperatore

$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = $changeDate
$ws2.Range("D3").Style = "Normal"
